$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 and 48: swap Polygon and InjectiveProtocol entries, update price/volume
$ws.Range("B47").Value = "Polygon"
$ws.Range("C47").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D47").Value = "0.411"
$ws.Range("E47").Value = "  +7.01%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "19.30"
$ws.Range("E48").Value = "  +1.60%  "

$ws.Range("D2").Value = "59.955.24"
$ws.Range("E2").Value = "  +2.94%  "
$ws.Range("D3").Value = "2.418.38"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("D5").Value = "551.72"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "137.20"
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.69%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  +3.78%  "
$ws.Range("D14").Value = "2.840.78"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "59.894.38"
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("D16").Value = "0.0000138"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "2.413.44"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").Value = "11.30"
$ws.Range("E18").Value = "  +4.12%  "
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").Value = "330.85"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "65.86"
$ws.Range("E23").Value = "  +3.65%  "
$ws.Range("D24").Value = "0.172"
$ws.Range("E24").Value = "  +3.73%  "
$ws.Range("E25").Value = "  +3.89%  "
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("E28").Value = "  +5.54%  "
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("D30").Value = "170.69"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "6.19"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +4.07%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("D41").Value = "314.54"
$ws.Range("E41").Value = "  +8.92%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "138.36"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "0.579"
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").Value = "17.63"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "11.05"
$ws.Range("E51").Value = "  -0.25%  "
